# turn 2 skill to spike skill #28
# Adds two new skill rows (57000013 "大墓地"/grave, 57000014 "德鲁伊"/sidekickall)
# to the "表1" table on the "标准" sheet, right below the existing 12 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: 57000013 / 大墓地 / grave -------------------------------------
$ws.Cells.Item(16, 1).Value = 57000013
$ws.Cells.Item(16, 2).Value = "大墓地"
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = "grave"
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = "'false"
$ws.Cells.Item(16, 9).Value = "'false"
$ws.Cells.Item(16, 10).Value = "'false"

# --- Row 17: 57000014 / 德鲁伊 / sidekickall --------------------------------
$ws.Cells.Item(17, 1).Value = 57000014
$ws.Cells.Item(17, 2).Value = "德鲁伊"
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = "sidekickall"
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = "'false"
$ws.Cells.Item(17, 9).Value = "'false"
$ws.Cells.Item(17, 10).Value = "'false"

# Match the formatting of the last existing data row (15) on the new rows,
# pasting formats only so the values just written above are kept intact.
$ws.Range("A15:J15").Copy()
$ws.Range("A16:J16").PasteSpecial(-4122)
$ws.Range("A15:J15").Copy()
$ws.Range("A17:J17").PasteSpecial(-4122)

# Grow the "表1" table so it covers the two new rows.
$lo = $ws.ListObjects.Item(1)
$null = $lo.Resize($ws.Range("A3:J17"))

# Match the author's final selection (cell F9) recorded in the saved file.
$null = $ws.Range("F9").Select()
